$wb = $excel.ActiveWorkbook

# --- Sheet: Matriz_Resultados ---
$ws1 = $wb.Worksheets.Item("Matriz_Resultados")
$ws1.Range("C2").Value = 0
$ws1.Range("D2").Value = 0
$ws1.Range("G2").Value = 0
$ws1.Range("H2").Value = 0
$ws1.Range("I2").Value = 0
$ws1.Range("J2").Value = 0
$ws1.Range("B3").Value = 0
$ws1.Range("I3").Value = 0
$ws1.Range("B4").Value = 0
$ws1.Range("E4").Value = 0
$ws1.Range("F4").Value = 0
$ws1.Range("D5").Value = 0
$ws1.Range("G5").Value = 0
$ws1.Range("H5").Value = 0
$ws1.Range("D6").Value = 0
$ws1.Range("G6").Value = 0
$ws1.Range("H6").Value = 0
$ws1.Range("J6").Value = 0
$ws1.Range("B7").Value = 0
$ws1.Range("E7").Value = 0
$ws1.Range("F7").Value = 0
$ws1.Range("B8").Value = 0
$ws1.Range("E8").Value = 0
$ws1.Range("F8").Value = 0
$ws1.Range("B9").Value = 0
$ws1.Range("C9").Value = 0
$ws1.Range("B10").Value = 0
$ws1.Range("F10").Value = 0

# --- Sheet: P_valores ---
$ws2 = $wb.Worksheets.Item("P_valores")
$ws2.Range("C2").Value = 0.001732280406275599
$ws2.Range("D2").Value = 0.03553690120744424
$ws2.Range("E2").Value = 0.09692510740376958
$ws2.Range("F2").Value = 0.802539175392045
$ws2.Range("G2").Value = 0.01274359335693198
$ws2.Range("H2").Value = 0.01570862913530968
$ws2.Range("I2").Value = 0.002082149059181138
$ws2.Range("J2").Value = 0.007195763463765159
$ws2.Range("B3").Value = 0.001732280406275599
$ws2.Range("D3").Value = [double]"4.310463097567663E-10"
$ws2.Range("E3").Value = 0.0005105860238883153
$ws2.Range("F3").Value = 0.0006236394131142742
$ws2.Range("G3").Value = [double]"4.35984300173331E-06"
$ws2.Range("H3").Value = [double]"2.833045852801419E-06"
$ws2.Range("I3").Value = 0.004749877926248258
$ws2.Range("J3").Value = [double]"9.854793757035907E-06"
$ws2.Range("B4").Value = 0.03553690120744424
$ws2.Range("C4").Value = [double]"4.310463097567663E-10"
$ws2.Range("E4").Value = 0.006354680427024029
$ws2.Range("F4").Value = 0.01814263608366296
$ws2.Range("G4").Value = 0.548437929691767
$ws2.Range("H4").Value = 0.3035459084401453
$ws2.Range("I4").Value = [double]"2.798695875050328E-08"
$ws2.Range("J4").Value = 0.0001323567720408203
$ws2.Range("B5").Value = 0.09692510740376958
$ws2.Range("C5").Value = 0.0005105860238883153
$ws2.Range("D5").Value = 0.006354680427024029
$ws2.Range("F5").Value = 0.1214972155947127
$ws2.Range("G5").Value = 0.00344898439426955
$ws2.Range("H5").Value = 0.003475674279061369
$ws2.Range("I5").Value = 0.0007641766127153282
$ws2.Range("J5").Value = 0.2603074031094359
$ws2.Range("B6").Value = 0.802539175392045
$ws2.Range("C6").Value = 0.0006236394131142742
$ws2.Range("D6").Value = 0.01814263608366296
$ws2.Range("E6").Value = 0.1214972155947127
$ws2.Range("G6").Value = 0.005202714774447337
$ws2.Range("H6").Value = 0.006917648542805743
$ws2.Range("I6").Value = 0.0007142770770791618
$ws2.Range("J6").Value = 0.003844282014359601
$ws2.Range("B7").Value = 0.01274359335693198
$ws2.Range("C7").Value = [double]"4.35984300173331E-06"
$ws2.Range("D7").Value = 0.548437929691767
$ws2.Range("E7").Value = 0.00344898439426955
$ws2.Range("F7").Value = 0.005202714774447337
$ws2.Range("H7").Value = 0.5223405372207397
$ws2.Range("I7").Value = [double]"2.560026612030697E-06"
$ws2.Range("J7").Value = [double]"6.160927502341096E-05"
$ws2.Range("B8").Value = 0.01570862913530968
$ws2.Range("C8").Value = [double]"2.833045852801419E-06"
$ws2.Range("D8").Value = 0.3035459084401453
$ws2.Range("E8").Value = 0.003475674279061369
$ws2.Range("F8").Value = 0.006917648542805743
$ws2.Range("G8").Value = 0.5223405372207397
$ws2.Range("I8").Value = [double]"4.977215776857946E-06"
$ws2.Range("J8").Value = [double]"6.882446870237757E-05"
$ws2.Range("B9").Value = 0.002082149059181138
$ws2.Range("C9").Value = 0.004749877926248258
$ws2.Range("D9").Value = [double]"2.798695875050328E-08"
$ws2.Range("E9").Value = 0.0007641766127153282
$ws2.Range("F9").Value = 0.0007142770770791618
$ws2.Range("G9").Value = [double]"2.560026612030697E-06"
$ws2.Range("H9").Value = [double]"4.977215776857946E-06"
$ws2.Range("J9").Value = [double]"1.306705973824762E-05"
$ws2.Range("B10").Value = 0.007195763463765159
$ws2.Range("C10").Value = [double]"9.854793757035907E-06"
$ws2.Range("D10").Value = 0.0001323567720408203
$ws2.Range("E10").Value = 0.2603074031094359
$ws2.Range("F10").Value = 0.003844282014359601
$ws2.Range("G10").Value = [double]"6.160927502341096E-05"
$ws2.Range("H10").Value = [double]"6.882446870237757E-05"
$ws2.Range("I10").Value = [double]"1.306705973824762E-05"

# --- Sheet: Estadisticos_DM ---
$ws3 = $wb.Worksheets.Item("Estadisticos_DM")
$ws3.Range("C2").Value = 3.860224414494249
$ws3.Range("D2").Value = 2.326211033602827
$ws3.Range("E2").Value = -1.779156773411489
$ws3.Range("F2").Value = -0.2548570019888785
$ws3.Range("G2").Value = 2.854212488590707
$ws3.Range("H2").Value = 2.747932236517209
$ws3.Range("I2").Value = 3.76701748786005
$ws3.Range("J2").Value = -3.142666658164647
$ws3.Range("B3").Value = -3.860224414494249
$ws3.Range("D3").Value = -15.18939526811226
$ws3.Range("E3").Value = -4.488212476949606
$ws3.Range("F3").Value = -4.38408058842305
$ws3.Range("G3").Value = -7.229568207045736
$ws3.Range("H3").Value = -7.510577134412886
$ws3.Range("I3").Value = -3.351488462820882
$ws3.Range("J3").Value = -6.715473455456671
$ws3.Range("B4").Value = -2.326211033602827
$ws3.Range("C4").Value = 15.18939526811226
$ws3.Range("E4").Value = -3.205184383521898
$ws3.Range("F4").Value = -2.674414636898157
$ws3.Range("G4").Value = 0.6149741020444036
$ws3.Range("H4").Value = 1.068103446154648
$ws3.Range("I4").Value = 11.01085298148083
$ws3.Range("J4").Value = -5.209135733142663
$ws3.Range("B5").Value = 1.779156773411489
$ws3.Range("C5").Value = 4.488212476949606
$ws3.Range("D5").Value = 3.205184383521898
$ws3.Range("F5").Value = 1.648516355217216
$ws3.Range("G5").Value = 3.512425980627033
$ws3.Range("H5").Value = 3.508546505851935
$ws3.Range("I5").Value = 4.278864438591858
$ws3.Range("J5").Value = -1.173155710531135
$ws3.Range("B6").Value = 0.2548570019888785
$ws3.Range("C6").Value = 4.38408058842305
$ws3.Range("D6").Value = 2.674414636898157
$ws3.Range("E6").Value = -1.648516355217216
$ws3.Range("G6").Value = 3.30572185712368
$ws3.Range("H6").Value = 3.162496044984905
$ws3.Range("I6").Value = 4.313762460547967
$ws3.Range("J6").Value = -3.457835655201027
$ws3.Range("B7").Value = -2.854212488590707
$ws3.Range("C7").Value = 7.229568207045736
$ws3.Range("D7").Value = -0.6149741020444036
$ws3.Range("E7").Value = -3.512425980627033
$ws3.Range("F7").Value = -3.30572185712368
$ws3.Range("H7").Value = 0.6561912319224935
$ws3.Range("I7").Value = 7.57760939856421
$ws3.Range("J7").Value = -5.634249118029534
$ws3.Range("B8").Value = -2.747932236517209
$ws3.Range("C8").Value = 7.510577134412886
$ws3.Range("D8").Value = -1.068103446154648
$ws3.Range("E8").Value = -3.508546505851935
$ws3.Range("F8").Value = -3.162496044984905
$ws3.Range("G8").Value = -0.6561912319224935
$ws3.Range("I8").Value = 7.144556562521533
$ws3.Range("J8").Value = -5.571827243408594
$ws3.Range("B9").Value = -3.76701748786005
$ws3.Range("C9").Value = 3.351488462820882
$ws3.Range("D9").Value = -11.01085298148083
$ws3.Range("E9").Value = -4.278864438591858
$ws3.Range("F9").Value = -4.313762460547967
$ws3.Range("G9").Value = -7.57760939856421
$ws3.Range("H9").Value = -7.144556562521533
$ws3.Range("J9").Value = -6.542655655293229
$ws3.Range("B10").Value = 3.142666658164647
$ws3.Range("C10").Value = 6.715473455456671
$ws3.Range("D10").Value = 5.209135733142663
$ws3.Range("E10").Value = 1.173155710531135
$ws3.Range("F10").Value = 3.457835655201027
$ws3.Range("G10").Value = 5.634249118029534
$ws3.Range("H10").Value = 5.571827243408594
$ws3.Range("I10").Value = 6.542655655293229

# --- Sheet: Resumen ---
$ws4 = $wb.Worksheets.Item("Resumen")
$ws4.Range("A2").Value = "Sieve Bootstrap"
$ws4.Range("B2").Value = 6
$ws4.Range("C2").Value = 0
$ws4.Range("D2").Value = 2
$ws4.Range("E2").Value = 75
$ws4.Range("F2").Value = 0.5464135857990908
$ws4.Range("A3").Value = "DeepAR"
$ws4.Range("B3").Value = 6
$ws4.Range("C3").Value = 0
$ws4.Range("D3").Value = 2
$ws4.Range("E3").Value = 75
$ws4.Range("F3").Value = 0.574004851152303
$ws4.Range("A4").Value = "AV-MCPS"
$ws4.Range("B4").Value = 1
$ws4.Range("C4").Value = 2
$ws4.Range("D4").Value = 5
$ws4.Range("E4").Value = 12.5
$ws4.Range("F4").Value = 0.6985225186644016
$ws4.Range("A5").Value = "MCPS"
$ws4.Range("B5").Value = 1
$ws4.Range("C5").Value = 2
$ws4.Range("D5").Value = 5
$ws4.Range("E5").Value = 12.5
$ws4.Range("F5").Value = 0.7059666671819314
$ws4.Range("A6").Value = "LSPM"
$ws4.Range("B6").Value = 1
$ws4.Range("C6").Value = 2
$ws4.Range("D6").Value = 5
$ws4.Range("E6").Value = 12.5
$ws4.Range("F6").Value = 0.7196613918393555
$ws4.Range("A7").Value = "Block Bootstrapping"
$ws4.Range("B7").Value = 0
$ws4.Range("C7").Value = 0
$ws4.Range("D7").Value = 8
$ws4.Range("E7").Value = 0
$ws4.Range("F7").Value = 0.9559532634470479
$ws4.Range("A8").Value = "LSPMW"
$ws4.Range("B8").Value = 0
$ws4.Range("C8").Value = 2
$ws4.Range("D8").Value = 6
$ws4.Range("E8").Value = 0
$ws4.Range("F8").Value = 1.076995155575413
$ws4.Range("A9").Value = "AREPD"
$ws4.Range("B9").Value = 0
$ws4.Range("C9").Value = 2
$ws4.Range("D9").Value = 6
$ws4.Range("E9").Value = 0
$ws4.Range("F9").Value = 0.9604152130268676
$ws4.Range("A10").Value = "EnCQR-LSTM"
$ws4.Range("B10").Value = 0
$ws4.Range("C10").Value = 5
$ws4.Range("D10").Value = 3
$ws4.Range("E10").Value = 0
$ws4.Range("F10").Value = 1.15301740025674
